$d = $word.ActiveDocument

# 1. Remove the old _GoBack bookmark (currently sits after "MP73010" in paragraph 1)
$d.Bookmarks("_GoBack").Delete()

# 2. Rewrite paragraph 4 (the ">>> ... >>>" line) with the new text and move
#    the _GoBack bookmark to sit between "Changes Made!!" and the closing ">>>"
$p = $d.Paragraphs(4)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = ""
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>&gt;&gt;&gt;</w:t></w:r><w:r><w:t>Changes Made!!</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>&gt;&gt;&gt;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

Write-Output ("done; paragraph4 text=[" + $p.Range.Text + "]")
